# Applies the "DTR summary cleanup + legends" edit:
#  1. Clears the stray "excess column" overtime-hour values in column I
#     for the rows where they shouldn't have been populated.
#  2. Adds a "Legends:" heading plus three colour-coded legend entries
#     (matching the existing highlight colours used in the DTR table)
#     below the summary block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Remove excess NO. OF OVERTIME HOURS values (column I) that were
#    left over in the summary rows.
# ---------------------------------------------------------------------
foreach ($cellRef in @("I7", "I8", "I9", "I10", "I14", "I15", "I16", "I17", "I18")) {
    $ws.Range($cellRef).ClearContents()
}

# ---------------------------------------------------------------------
# 2. "Legends:" title, formatted like the report's other headings
#    (re-uses the same large bold-underlined look as A1:A3).
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("E24").PasteSpecial(-4122) | Out-Null
$ws.Range("E24").Value = "Legends:"
$ws.Range("E24:P24").Merge()

# ---------------------------------------------------------------------
# 3. Legend entry 1 (cyan) - request/remark note
# ---------------------------------------------------------------------
$ws.Range("E25").Interior.Color = 13411113
$ws.Range("E25:E26").Merge()

$ws.Range("F25").Value = "Employee has request(s)/remark(s) for that day.`n*May incur late and/or undertime depending on his or her time-in and time-out."
$ws.Range("F25:P26").Font.Bold = $true
$ws.Range("F25:P26").Font.Underline = $true
$ws.Range("F25:P26").Merge()

# ---------------------------------------------------------------------
# 4. Legend entry 2 (orange) - half-day note
# ---------------------------------------------------------------------
$ws.Range("E27").Interior.Color = 6737151
$ws.Range("E27:E28").Merge()

$ws.Range("F27").Value = "Employee is considered half-day because of his time-in or time-out."
$ws.Range("F27:P28").Font.Bold = $true
$ws.Range("F27:P28").Font.Underline = $true
$ws.Range("F27:P28").Merge()

# ---------------------------------------------------------------------
# 5. Legend entry 3 (red) - absent note
# ---------------------------------------------------------------------
$ws.Range("E29").Interior.Color = 6184671
$ws.Range("E29:E30").Merge()

$ws.Range("F29").Value = "Employee has no time-in and therefore, considered as absent."
$ws.Range("F29:P30").Font.Bold = $true
$ws.Range("F29:P30").Font.Underline = $true
$ws.Range("F29:P30").Merge()
